# "Generate Report for Handback"
# The en-US source has been handed back to zh-cn and de-de: the overview
# status moves from "Ready for handoff" to "Handed back: in sync with en-US",
# and the per-locale sheets get their Latest Target File / Latest Handback
# File / Latest Handback DateTime columns filled in (with a hyperlink on the
# target-file cell), plus the corresponding columns are widened.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"
$mdFile = "1442afcb-9d0e-4476-99c2-45f1560af881.md"
$mdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d60bf023c839e1118d2a2ea522c73ca4f4ccd1d6/e2e/1442afcb-9d0e-4476-99c2-45f1560af881.md"

# Excel's ColumnWidth property is expressed in characters; the stored XML
# column width is ColumnWidth + 5/6. Pick the ColumnWidth that reproduces
# the desired stored widths.
$wideWidth = 29.9777047293527 - (5/6)
$fortyWidth = 40 - (5/6)

# ---------------------------------------------------------------------
# Overview sheet: status text for the zh-cn / de-de columns, and widen
# those two columns to fit the longer text.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value2 = $newStatus
$wsOverview.Range("F2").Value2 = $newStatus
$wsOverview.Range("E1").ColumnWidth = $wideWidth
$wsOverview.Range("F1").ColumnWidth = $wideWidth

# ---------------------------------------------------------------------
# zh-cn sheet: mark handed back, fill in target/handback file + datetime.
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value2 = $newStatus
$wsZhCn.Range("J2").Value2 = "1442afcb-9d0e-4476-99c2-45f1560af881.88c59249ca584fd795a79250abd1c6345bfec02d.zh-cn.xlf"
$wsZhCn.Range("K2").Value2 = "2016-09-03 19:12:37"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $mdUrl, "", "", $mdFile) | Out-Null

$wsZhCn.Range("C1").ColumnWidth = $wideWidth
$wsZhCn.Range("I1").ColumnWidth = $fortyWidth
$wsZhCn.Range("J1").ColumnWidth = $fortyWidth

# ---------------------------------------------------------------------
# de-de sheet: mark handed back, fill in target/handback file + datetime.
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value2 = $newStatus
$wsDeDe.Range("J2").Value2 = "1442afcb-9d0e-4476-99c2-45f1560af881.88c59249ca584fd795a79250abd1c6345bfec02d.de-de.xlf"
$wsDeDe.Range("K2").Value2 = "2016-09-03 19:12:44"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $mdUrl, "", "", $mdFile) | Out-Null

$wsDeDe.Range("C1").ColumnWidth = $wideWidth
$wsDeDe.Range("I1").ColumnWidth = $fortyWidth
$wsDeDe.Range("J1").ColumnWidth = $fortyWidth
